# Update latest output (run 128)
$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# ---- Schedule sheet (row 3-5) ----
$wsSchedule.Range("E3").Value = -30.35730749999998
$wsSchedule.Range("F3").Value = -0.8923370811287473

$wsSchedule.Range("B4").Value = 46043.14583333334
$wsSchedule.Range("C4").Value = 5.5
$wsSchedule.Range("D4").Value = 20.79
$wsSchedule.Range("E4").Value = 619.4632859999999
$wsSchedule.Range("F4").Value = 29.79621385281385

$wsSchedule.Range("A5").Value = 46043.3125
$wsSchedule.Range("C5").Value = 8.5
$wsSchedule.Range("D5").Value = 32.13
$wsSchedule.Range("E5").Value = -226.12220475
$wsSchedule.Range("F5").Value = -7.037728127917834

# ---- Detailed sheet ----
$wsDetailed.Range("B33").Value = -6.39804
$wsDetailed.Range("B34").Value = -5.51

$wsDetailed.Range("B35").Value = 22.07
$wsDetailed.Range("C35").Value = "historical"

$wsDetailed.Range("B36").Value = 36.07
$wsDetailed.Range("C36").Value = "historical"

$wsDetailed.Range("B37").Value = 50.8425
$wsDetailed.Range("C37").Value = "historical"

$wsDetailed.Range("B38").Value = 60.35108
$wsDetailed.Range("B39").Value = 63.74797
$wsDetailed.Range("B40").Value = 73.94592
$wsDetailed.Range("B41").Value = 75.72364
$wsDetailed.Range("B42").Value = 78
$wsDetailed.Range("B43").Value = 76.20085
$wsDetailed.Range("B44").Value = 69.69265
$wsDetailed.Range("B45").Value = 63.79708
$wsDetailed.Range("B47").Value = 64.26672000000001
$wsDetailed.Range("B48").Value = 57.3
$wsDetailed.Range("B49").Value = 57.06003

$wsDetailed.Range("E56").Value = "ON"

$wsDetailed.Range("B58").Value = 65.94062
$wsDetailed.Range("B59").Value = 67.70805
$wsDetailed.Range("B60").Value = 72.01038
$wsDetailed.Range("B61").Value = 79.36201
$wsDetailed.Range("B62").Value = 80.45009
$wsDetailed.Range("B63").Value = 68.98878999999999

$wsDetailed.Range("B64").Value = 57.06012
$wsDetailed.Range("E64").Value = "OFF"

$wsDetailed.Range("B65").Value = 0.7
$wsDetailed.Range("B66").Value = 0.009719999999999999

$wsDetailed.Range("B68").Value = -6.62812
$wsDetailed.Range("B69").Value = -6.99226
$wsDetailed.Range("B70").Value = -8.691649999999999
$wsDetailed.Range("B71").Value = -8.85641
$wsDetailed.Range("B72").Value = -14
$wsDetailed.Range("B73").Value = -14.49854
$wsDetailed.Range("B74").Value = -16.35771
$wsDetailed.Range("B75").Value = -23.07171
$wsDetailed.Range("B76").Value = -22.40654

$wsDetailed.Range("B78").Value = -23.5
$wsDetailed.Range("B79").Value = -23.78763
$wsDetailed.Range("B80").Value = -23.5
$wsDetailed.Range("B81").Value = -14
$wsDetailed.Range("B82").Value = -6.40066
$wsDetailed.Range("B83").Value = -4.61188

$wsDetailed.Range("B85").Value = 46.9682
$wsDetailed.Range("B86").Value = 56.52321
$wsDetailed.Range("B87").Value = 65
$wsDetailed.Range("B88").Value = 77.94
$wsDetailed.Range("B89").Value = 105.79
$wsDetailed.Range("B90").Value = 79.95
$wsDetailed.Range("B91").Value = 73.82011
$wsDetailed.Range("B92").Value = 70.12006

$wsDetailed.Range("B94").Value = 60.27702
